# GoldDrybarTestData.xlsx edit: insert a new reward-tier row into the
# "DataSet" sheet at row 36 ("$20 Off (400 points)"), which pushes every
# row from 36..120 down by one (->37..121). The sheet's own Insert() call
# shifts cell data/styles correctly, but this engine's Hyperlinks
# collection is NOT shifted automatically by a row insert, so the
# worksheet-level hyperlinks are rebuilt from scratch afterwards with the
# correct (possibly shifted) target ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSet")

# ---------------------------------------------------------------------
# 1) Capture the hyperlink target (mailto:/https:) for every existing
#    hyperlink together with its NEW (post-insert) destination range.
#    Rows >= 36 shift down by one; everything above row 36 is unchanged.
#
#    NOTE: this engine's Hyperlinks.Add(... TextToDisplay) clobbers the
#    *cell's own value* with the display text instead of only annotating
#    the hyperlink (verified experimentally: passing a "display" string
#    rewrites <c> to a shared-string pointing at that display text). Since
#    several of the recorded display strings differ from the real cell
#    text (e.g. a "mailto:" prefix that isn't part of the cell value), the
#    display text is intentionally left unset everywhere so the already
#    shifted-down cell contents are never overwritten.
# ---------------------------------------------------------------------
$hyperlinkSpecs = @(
    @{Ref='P23'; Target='mailto:avayugundla@helenoftroy.com'}
    @{Ref='P24'; Target='mailto:avayugundla@helenoftroy.com'}
    @{Ref='E2'; Target='mailto:Lotuswave@123'}
    @{Ref='P31'; Target='mailto:avayugundla@helenoftroy.com'}
    @{Ref='P37'; Target='mailto:Paypal-buyer@hydroflask.com'}
    @{Ref='P39'; Target='mailto:skasarla@helenoftroy.com'}
    @{Ref='P25'; Target='mailto:avayugundla@helenoftroy.com'}
    @{Ref='P53'; Target='mailto:qatesting.lotuswave@gmail.com'}
    @{Ref='P26:P27'; Target='mailto:avayugundla@helenoftroy.com'}
    @{Ref='F10'; Target='mailto:!#@'}
    @{Ref='D10'; Target='mailto:Lotus@1'}
    @{Ref='E10'; Target='mailto:Lotus@1235'}
    @{Ref='B56'; Target='mailto:hydroflaskemea978+7@gmail.com'}
    @{Ref='C56'; Target='mailto:hydroflaskemea978+7@gmail.com'}
    @{Ref='D56'; Target='mailto:Lotus@123'}
    @{Ref='E56'; Target='mailto:Lotus@123'}
    @{Ref='P56'; Target='mailto:hydroflaskemea978+7@gmail.com'}
    @{Ref='B57'; Target='mailto:hydroflaskemea978+8@gmail.com'}
    @{Ref='D57'; Target='mailto:Lotus@123'}
    @{Ref='E57'; Target='mailto:Lotuswave@1234'}
    @{Ref='P57'; Target='mailto:avayugundla@helenoftroy.com'}
    @{Ref='Q57'; Target='mailto:hydroflaskemea978+2@gmail.com'}
    @{Ref='B60'; Target='mailto:skasarla@helenoftroy.com'}
    @{Ref='B3'; Target='mailto:gsapram@helenoftroy.com'}
    @{Ref='C3'; Target='mailto:gsapram@helenoftroy.com'}
    @{Ref='D3'; Target='mailto:Lotuswave@123'}
    @{Ref='E3'; Target='mailto:Lotuswave@123'}
    @{Ref='B2'; Target='mailto:avayugundla@helenoftroy.com'}
    @{Ref='C2'; Target='mailto:avayugundla@helenoftroy.com'}
    @{Ref='D2'; Target='mailto:Lotuswave@123'}
    @{Ref='AE8'; Target='https://mcloud-na-stage4.drybar.com/the-straight-shot-blow-drying-flat-iron.html'}
    @{Ref='AO10'; Target='https://mcloud-na-stage4.drybar.com/hydroflask'}
    @{Ref='E54'; Target='mailto:Lotuswave@123'}
    @{Ref='B54'; Target='mailto:avayugundla@helenoftroy.com'}
    @{Ref='C54'; Target='mailto:avayugundla@helenoftroy.com'}
    @{Ref='D54'; Target='mailto:Lotuswave@123'}
    @{Ref='P61'; Target='mailto:avayugundla@helenoftroy.com'}
    @{Ref='P62'; Target='mailto:avayugundla@helenoftroy.com'}
    @{Ref='B52'; Target='mailto:avayugundla@helenoftroy.com'}
    @{Ref='C52'; Target='mailto:avayugundla@helenoftroy.com'}
    @{Ref='D52'; Target='mailto:Lotuswave@123'}
)

# ---------------------------------------------------------------------
# 2) Drop all existing hyperlinks up front -- the underlying range/row
#    data they point at is about to move, and this engine only knows how
#    to clear the whole collection at once (a per-hyperlink Delete() is a
#    no-op here).
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------
# 3) Insert the new row. Excel shifts row 36..120 down to 37..121 and
#    carries the per-cell formatting along for the ride.
# ---------------------------------------------------------------------
$ws.Rows.Item(36).Insert()

# ---------------------------------------------------------------------
# 4) Populate the newly-inserted row 36 with the new reward tier.
# ---------------------------------------------------------------------
$ws.Range("A36").Value = "`$20 Off (400 points)"
$ws.Range("AL36").Value = 400
$ws.Range("AM36").Value = "points"

# ---------------------------------------------------------------------
# 5) Recreate every hyperlink at its (possibly shifted) destination.
# ---------------------------------------------------------------------
foreach ($spec in $hyperlinkSpecs) {
    $target = $ws.Range($spec.Ref)
    $ws.Hyperlinks.Add($target, $spec.Target)
}

# ---------------------------------------------------------------------
# 6) Match the recorded selection state for the sheet.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("AM42").Select()
